$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (old row 19) - the table shrinks from 18 to 17 data rows
$ws.Rows.Item(19).Delete()

# Rewrite the data rows (A2:C18) with the new player/position/team data
$data = @(
    @("De'Aaron Fox",   "PG",          "Sacramento Kings"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Tyler Herro",    "PG,SG",       "Miami Heat"),
    @("Ja Morant",      "PG",          "Memphis Grizzlies"),
    @("DeMar DeRozan",  "SF,PF",       "Sacramento Kings"),
    @("Miles Bridges",  "SF,PF",       "Charlotte Hornets"),
    @("Josh Giddey",    "PG,SG,SF",    "Chicago Bulls"),
    @("Nikola Vucevic", "PF,C",        "Chicago Bulls"),
    @("Zach Edey",      "C",           "Memphis Grizzlies"),
    @("Nick Richards",  "C",           "Phoenix Suns"),
    @("Brook Lopez",    "C",           "Milwaukee Bucks"),
    @("Clint Capela",   "C",           "Atlanta Hawks"),
    @("Shaedon Sharpe", "SG,SF",       "Portland Trail Blazers"),
    @("Mikal Bridges",  "SG,SF,PF",    "New York Knicks"),
    @("Luka Doncic",    "PG,SG",       "Dallas Mavericks"),
    @("Evan Mobley",    "PF,C",        "Cleveland Cavaliers"),
    @("Bobby Portis",   "PF,C",        "Milwaukee Bucks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
